# Apply "Natmi following Dr Hou advice" data update
# Expand Icam1-Spn LR-pair table from 4 data rows (target cluster = M2 only)
# to 8 data rows (target cluster = M2 and sCs, for each of the 4 sending clusters),
# and refresh all computed metric columns (E..T) with the revised values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Icam1"
$ws.Range("C2").Value2 = "Spn"
$ws.Range("D2").Value2 = "M2"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 14.452709
$ws.Range("H2").Value2 = 43.358127
$ws.Range("I2").Value2 = 0.1476906377370901
$ws.Range("J2").Value2 = 0.1476906377370901
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.400925
$ws.Range("N2").Value2 = 7.202775
$ws.Range("O2").Value2 = 0.9592478079643895
$ws.Range("P2").Value2 = 0.9592478079643894
$ws.Range("Q2").Value2 = 34.699870355825
$ws.Range("R2").Value2 = 312.298833202425
$ws.Range("S2").Value2 = 0.1416719205061664
$ws.Range("T2").Value2 = 0.1416719205061664

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Icam1"
$ws.Range("C3").Value2 = "Spn"
$ws.Range("D3").Value2 = "sCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 14.452709
$ws.Range("H3").Value2 = 43.358127
$ws.Range("I3").Value2 = 0.1476906377370901
$ws.Range("J3").Value2 = 0.1476906377370901
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.1019996666666667
$ws.Range("N3").Value2 = 0.305999
$ws.Range("O3").Value2 = 0.04075219203561061
$ws.Range("P3").Value2 = 0.04075219203561061
$ws.Range("Q3").Value2 = 1.474171500430333
$ws.Range("R3").Value2 = 13.267543503873
$ws.Range("S3").Value2 = 0.006018717230923695
$ws.Range("T3").Value2 = 0.006018717230923696

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Icam1"
$ws.Range("C4").Value2 = "Spn"
$ws.Range("D4").Value2 = "M2"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 45.91529066666667
$ws.Range("H4").Value2 = 137.745872
$ws.Range("I4").Value2 = 0.4692032864180593
$ws.Range("J4").Value2 = 0.4692032864180593
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 2.400925
$ws.Range("N4").Value2 = 7.202775
$ws.Range("O4").Value2 = 0.9592478079643895
$ws.Range("P4").Value2 = 0.9592478079643894
$ws.Range("Q4").Value2 = 110.2391692438667
$ws.Range("R4").Value2 = 992.1525231948001
$ws.Range("S4").Value2 = 0.4500822239862109
$ws.Range("T4").Value2 = 0.4500822239862109

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Icam1"
$ws.Range("C5").Value2 = "Spn"
$ws.Range("D5").Value2 = "sCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 45.91529066666667
$ws.Range("H5").Value2 = 137.745872
$ws.Range("I5").Value2 = 0.4692032864180593
$ws.Range("J5").Value2 = 0.4692032864180593
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.1019996666666667
$ws.Range("N5").Value2 = 0.305999
$ws.Range("O5").Value2 = 0.04075219203561061
$ws.Range("P5").Value2 = 0.04075219203561061
$ws.Range("Q5").Value2 = 4.683344342903112
$ws.Range("R5").Value2 = 42.15009908612801
$ws.Range("S5").Value2 = 0.01912106243184836
$ws.Range("T5").Value2 = 0.01912106243184836

# Row 6
$ws.Range("A6").Value2 = "M2"
$ws.Range("B6").Value2 = "Icam1"
$ws.Range("C6").Value2 = "Spn"
$ws.Range("D6").Value2 = "M2"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 36.015061
$ws.Range("H6").Value2 = 108.045183
$ws.Range("I6").Value2 = 0.3680339324088102
$ws.Range("J6").Value2 = 0.3680339324088103
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.400925
$ws.Range("N6").Value2 = 7.202775
$ws.Range("O6").Value2 = 0.9592478079643895
$ws.Range("P6").Value2 = 0.9592478079643894
$ws.Range("Q6").Value2 = 86.46946033142501
$ws.Range("R6").Value2 = 778.2251429828251
$ws.Range("S6").Value2 = 0.3530357429196655
$ws.Range("T6").Value2 = 0.3530357429196655

# Row 7
$ws.Range("A7").Value2 = "M2"
$ws.Range("B7").Value2 = "Icam1"
$ws.Range("C7").Value2 = "Spn"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 36.015061
$ws.Range("H7").Value2 = 108.045183
$ws.Range("I7").Value2 = 0.3680339324088102
$ws.Range("J7").Value2 = 0.3680339324088103
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.1019996666666667
$ws.Range("N7").Value2 = 0.305999
$ws.Range("O7").Value2 = 0.04075219203561061
$ws.Range("P7").Value2 = 0.04075219203561061
$ws.Range("Q7").Value2 = 3.673524216979667
$ws.Range("R7").Value2 = 33.06171795281701
$ws.Range("S7").Value2 = 0.01499818948914477
$ws.Range("T7").Value2 = 0.01499818948914477

# Row 8
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Icam1"
$ws.Range("C8").Value2 = "Spn"
$ws.Range("D8").Value2 = "M2"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 1.474929666666667
$ws.Range("H8").Value2 = 4.424789000000001
$ws.Range("I8").Value2 = 0.01507214343604052
$ws.Range("J8").Value2 = 0.01507214343604052
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 2.400925
$ws.Range("N8").Value2 = 7.202775
$ws.Range("O8").Value2 = 0.9592478079643895
$ws.Range("P8").Value2 = 0.9592478079643894
$ws.Range("Q8").Value2 = 3.541195509941667
$ws.Range("R8").Value2 = 31.870759589475
$ws.Range("S8").Value2 = 0.01445792055234673
$ws.Range("T8").Value2 = 0.01445792055234673

# Row 9
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Icam1"
$ws.Range("C9").Value2 = "Spn"
$ws.Range("D9").Value2 = "sCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 1.474929666666667
$ws.Range("H9").Value2 = 4.424789000000001
$ws.Range("I9").Value2 = 0.01507214343604052
$ws.Range("J9").Value2 = 0.01507214343604052
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.1019996666666667
$ws.Range("N9").Value2 = 0.305999
$ws.Range("O9").Value2 = 0.04075219203561061
$ws.Range("P9").Value2 = 0.04075219203561061
$ws.Range("Q9").Value2 = 0.1504423343567778
$ws.Range("R9").Value2 = 1.353981009211
$ws.Range("S9").Value2 = 0.0006142228836937913
$ws.Range("T9").Value2 = 0.0006142228836937913
